$d = $word.ActiveDocument

# --- 1. Entity attribute bullet list: shift items down, inserting "Работник" ---
# Original order: Лаборатория, Человек, Действие, Вид исследований, Вид воображения, Должность
# Target order:   Лаборатория, Человек, Работник, Действие, Вид исследований, Вид воображения, Должность

# Rename existing bullets in place (first occurrence only - the later
# "Ассоциативные"/"Характеристические" sentences reuse these words too).
$d.Content.Find.Execute("Действие", $true, $false, $false, $false, $false, $true, 1, $false, "Работник", 1) | Out-Null
$d.Content.Find.Execute("Вид исследований", $true, $false, $false, $false, $false, $true, 1, $false, "Действие", 1) | Out-Null
$d.Content.Find.Execute("Вид воображения", $true, $false, $false, $false, $false, $true, 1, $false, "Вид исследований", 1) | Out-Null

# Insert a brand-new bullet ("Вид воображения") right after the bullet that
# now reads "Вид исследований" (formerly "Вид воображения"), before "Должность".
# (Split into two InsertAfter calls - a single call with an embedded CR
# misplaces the text that follows the break.)
$rng = $d.Content
$rng.Find.Execute("Вид исследований") | Out-Null
$rng.Collapse(0)
$rng.InsertAfter("`r")
$rng2 = $d.Range($rng.End, $rng.End)
$rng2.InsertAfter("Вид воображения")

Write-Output "list-reorder-done"

# --- 2. "Ассоциативные" relationship bullet gains a second clause about
#        the "Работник" association. ---
$rng = $d.Content
$rng.Find.Execute("Ассоциативные – Действие (связывает Лабораторию и Человека)") | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(", Работник (связывает Человека и Лабораторию)")

Write-Output "associative-bullet-done"
